$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("D1")
    if ($cell.Text -eq "tot_fronds") {
        $cell.Value = "total_fronds"
    }
}
